# Applies the "Automatic update of files" edit to the TROSA overview sheet:
#   1. Column C ("Förändrad") date serial is bumped from 45184 to 45186 for
#      every data row.
#   2. Every HYPERLINK() formula in columns S-Y gets a second argument added
#      (the "friendly name" text), equal to that row's column A value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 135
$hyperlinkCols = @(19, 20, 21, 22, 23, 24, 25)  # S, T, U, V, W, X, Y

for ($r = $firstRow; $r -le $lastRow; $r++) {

    # --- 1. Bump the "Förändrad" date in column C ---
    $cCell = $ws.Cells.Item($r, 3)
    if ($cCell.Value2 -eq 45184) {
        $cCell.Value = 45186
    }

    # --- 2. Add the friendly-name argument to HYPERLINK formulas ---
    $label = $ws.Cells.Item($r, 1).Value2

    foreach ($col in $hyperlinkCols) {
        $cell = $ws.Cells.Item($r, $col)
        $formula = $cell.Formula
        if ($formula -and $formula.StartsWith("=HYPERLINK(") -and -not $formula.Contains(",")) {
            $newFormula = $formula.Substring(0, $formula.Length - 1) + ', "' + $label + '")'
            $cell.Formula = $newFormula
        }
    }
}
